$wb = $excel.ActiveWorkbook

# --- 1. Flip test.runparallel value from true to false on c-demo_ui sheet ---
$demoSheet = $wb.Worksheets.Item("c-demo_ui")
$demoSheet.Cells.Item(29, 2).Value = "'false"

# --- 2. Add the "t-exceltoexcel" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet5 = $wb.Worksheets.Add($null, $lastSheet)
$sheet5.Name = "t-exceltoexcel"

$sheet5.Cells.Item(1, 1).Value = "Step"
$sheet5.Cells.Item(1, 2).Value = "Action"
$sheet5.Cells.Item(1, 3).Value = "Target"
$sheet5.Cells.Item(1, 4).Value = "Input"
$sheet5.Cells.Item(1, 5).Value = "Overrides"

$sheet5.Cells.Item(2, 1).Value = "config setup"
$sheet5.Cells.Item(2, 2).Value = "config"

$sheet5.Cells.Item(3, 1).Value = "execute component"
$sheet5.Cells.Item(3, 2).Value = "component"
$sheet5.Cells.Item(3, 3).Value = "test.soaptest"

# --- 3. Add the "t-exceltocsv" sheet at the end ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet6 = $wb.Worksheets.Add($null, $lastSheet2)
$sheet6.Name = "t-exceltocsv"

$sheet6.Cells.Item(1, 1).Value = "Step"
$sheet6.Cells.Item(1, 2).Value = "Action"
$sheet6.Cells.Item(1, 3).Value = "Target"
$sheet6.Cells.Item(1, 4).Value = "Input"
$sheet6.Cells.Item(1, 5).Value = "Overrides"

$sheet6.Cells.Item(2, 1).Value = "config setup"
$sheet6.Cells.Item(2, 2).Value = "config"

$sheet6.Cells.Item(3, 1).Value = "execute component"
$sheet6.Cells.Item(3, 2).Value = "component"
$sheet6.Cells.Item(3, 3).Value = "test.alpha"

# --- 4. Make t-exceltocsv the active sheet ---
$sheet6.Activate()
